{"js": "// Replace each \"a\u00f7b=\" arithmetic prompt in the document's table cells with\n// its updated value, per the authored diff. Every \"before\" string is unique\n// within the document, so a scoped search-and-replace for each pair is\n// unambiguous and order-independent.\nconst replacements = [\n  [\"17\u00f76=\", \"96\u00f79=\"],\n  [\"76\u00f72=\", \"71\u00f73=\"],\n  [\"54\u00f74=\", \"15\u00f75=\"],\n  [\"48\u00f72=\", \"47\u00f72=\"],\n  [\"95\u00f79=\", \"72\u00f79=\"],\n  [\"54\u00f73=\", \"84\u00f74=\"],\n  [\"34\u00f76=\", \"87\u00f78=\"],\n  [\"84\u00f78=\", \"87\u00f77=\"],\n  [\"44\u00f73=\", \"75\u00f78=\"],\n  [\"93\u00f78=\", \"56\u00f78=\"],\n  [\"23\u00f73=\", \"57\u00f76=\"],\n  [\"70\u00f75=\", \"74\u00f76=\"],\n  [\"40\u00f72=\", \"51\u00f76=\"],\n  [\"42\u00f77=\", \"12\u00f75=\"],\n  [\"23\u00f76=\", \"71\u00f73=\"],\n  [\"38\u00f79=\", \"29\u00f76=\"],\n  [\"65\u00f74=\", \"90\u00f72=\"],\n  [\"21\u00f73=\", \"90\u00f77=\"],\n  [\"94\u00f73=\", \"32\u00f74=\"],\n  [\"88\u00f72=\", \"81\u00f73=\"],\n  [\"41\u00f73=\", \"29\u00f76=\"],\n  [\"31\u00f76=\", \"22\u00f77=\"],\n  [\"56\u00f73=\", \"26\u00f76=\"],\n  [\"75\u00f75=\", \"24\u00f76=\"],\n  [\"27\u00f74=\", \"96\u00f73=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Could not find text to replace: ${oldText}`);\n  }\n\n  results.items[0].insertText(newText, \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Replace each \"a\u00f7b=\" arithmetic prompt in the document's table cells with\n# its updated value, per the authored diff. Every \"before\" string is unique\n# within the document, so Find/Replace on the whole-document range for each\n# pair is unambiguous and order-independent.\n$d = $word.ActiveDocument\n\n$replacements = @(\n  @(\"17\u00f76=\", \"96\u00f79=\"),\n  @(\"76\u00f72=\", \"71\u00f73=\"),\n  @(\"54\u00f74=\", \"15\u00f75=\"),\n  @(\"48\u00f72=\", \"47\u00f72=\"),\n  @(\"95\u00f79=\", \"72\u00f79=\"),\n  @(\"54\u00f73=\", \"84\u00f74=\"),\n  @(\"34\u00f76=\", \"87\u00f78=\"),\n  @(\"84\u00f78=\", \"87\u00f77=\"),\n  @(\"44\u00f73=\", \"75\u00f78=\"),\n  @(\"93\u00f78=\", \"56\u00f78=\"),\n  @(\"23\u00f73=\", \"57\u00f76=\"),\n  @(\"70\u00f75=\", \"74\u00f76=\"),\n  @(\"40\u00f72=\", \"51\u00f76=\"),\n  @(\"42\u00f77=\", \"12\u00f75=\"),\n  @(\"23\u00f76=\", \"71\u00f73=\"),\n  @(\"38\u00f79=\", \"29\u00f76=\"),\n  @(\"65\u00f74=\", \"90\u00f72=\"),\n  @(\"21\u00f73=\", \"90\u00f77=\"),\n  @(\"94\u00f73=\", \"32\u00f74=\"),\n  @(\"88\u00f72=\", \"81\u00f73=\"),\n  @(\"41\u00f73=\", \"29\u00f76=\"),\n  @(\"31\u00f76=\", \"22\u00f77=\"),\n  @(\"56\u00f73=\", \"26\u00f76=\"),\n  @(\"75\u00f75=\", \"24\u00f76=\"),\n  @(\"27\u00f74=\", \"96\u00f73=\"),\n)\n\nforeach ($pair in $replacements) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $oldText\n  $find.Replacement.Text = $newText\n\n  $found = $find.Execute(\n    [ref]$oldText,   # FindText\n    $true,           # MatchCase\n    $false,          # MatchWholeWord\n    $false,          # MatchWildcards\n    $false,          # MatchSoundsLike\n    $false,          # MatchAllWordForms\n    $true,           # Forward\n    1,               # Wrap (wdFindContinue)\n    $false,          # Format\n    $newText,        # ReplaceWith\n    2                # Replace (wdReplaceAll)\n  )\n\n  if (-not $found) {\n    Write-Output \"WARNING: not found -> $oldText\"\n  }\n}\n\nWrite-Output \"done\"\n\n"}
